# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest scraped counts (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 651
    $ws.Range("F3").Value = 501
    $ws.Range("F4").Value = 36
    $ws.Range("F8").Value = 2046
    $ws.Range("F9").Value = 4097
    $ws.Range("F10").Value = 96
}
